$d = $word.ActiveDocument

# Helper: run a Find/Replace across the whole document content.
function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $old"
    }
}

$rsquo = [char]0x2019

# 1. Title: "MintTrack" + " Test Plan" runs merge into one run (text unchanged).
#    No textual change, nothing to do via Find/Replace (content identical).

# 2. Intro paragraph: "The " + "MintTrack" + " android application..." merge (text unchanged).

# 3. "applications memory footprint" -> "application's memory footprint"
Replace-Text "analyze the applications memory footprint" ("analyze the application" + $rsquo + "s memory footprint")

# 4. "framework incorrectly it rather likely" -> "framework incorrectly is rather likely"
Replace-Text "using the framework incorrectly it rather likely." "using the framework incorrectly is rather likely."

# 5. "are also display on the home tab" -> "are also displayed on the home tab"
Replace-Text "are also display on the home tab" "are also displayed on the home tab"

# 6. "deleteing" -> "deleting"
Replace-Text "buttons for deleteing and editing" "buttons for deleting and editing"

# 7. "a consistent way of access the database" -> "a consistent way of accessing the database"
Replace-Text "a consistent way of access the database" "a consistent way of accessing the database"

# 8. "Database Layer : Low level classes" -> "Database Layer: Low level classes"
Replace-Text "Database Layer : Low level classes" "Database Layer: Low level classes"

# 9. "These layer should be hidden" -> "These layers should be hidden"
Replace-Text "These layer should be hidden" "These layers should be hidden"

# 10. "Features to be Test" -> "Features to be Tested"
Replace-Text "Features to be Test" "Features to be Tested"

# 11. Expense paragraph corrections
Replace-Text "is considered that `"From Account`"" "is considered the `"From Account`""
Replace-Text "An expense is any type of financial event that would cost a person money, for example" ("An expense is any type of financial event that decreases a person" + $rsquo + "s money on hand, for example")

# 12. Transfer: "adding it to another account" -> "add it to another account"
Replace-Text "take currency from one account and adding it to another account" "take currency from one account and add it to another account"

# 13. Edit Transaction: "modify all data element, save" -> "modify all data elements, save"
Replace-Text "modify all data element, save any modifications" "modify all data elements, save any modifications"

# 14. Support Transactions: "all transaction created" -> "all transactions created"; "is display row" -> "is displayed row"
Replace-Text "supports the display of all transaction created via the entry tab" "supports the display of all transactions created via the entry tab"
Replace-Text "specific transaction type is display row by row per transaction" "specific transaction type is displayed row by row per transaction"

# 15. Transaction Interaction: add comma after "database"
Replace-Text "After querying the database if a transaction is pressed" "After querying the database, if a transaction is pressed"

# 16. Delete Action corrections
Replace-Text "The delete action will remove the exists of the transaction from the database" "The delete action will remove the existence of the transaction from the database"
Replace-Text "accordingly base on what is being deleted" "accordingly based on what is being deleted."

# 17. Scroll: "to many" -> "too many"; "are display to be displayed" -> "are to be displayed"
Replace-Text "When to many transactions are display to be displayed on the screen" "When too many transactions are to be displayed on the screen"

# 18. "Deactivating an category" -> "Deactivating a category"
Replace-Text "Deactivating an category will hide it from the category drop downs" "Deactivating a category will hide it from the category drop downs"
